$d = $word.ActiveDocument

[void]$d.Content.Find.Execute("2023-07-10 Monday", $false, $false, $false, $false, $false, $true, 1, $false, "2023-07-11 Tuesday", 2)
[void]$d.Content.Find.Execute("27-25=2", $false, $false, $false, $false, $false, $true, 1, $false, "74+25=99", 2)
[void]$d.Content.Find.Execute("31+19=50", $false, $false, $false, $false, $false, $true, 1, $false, "21+34=55", 2)
[void]$d.Content.Find.Execute("51+39=90", $false, $false, $false, $false, $false, $true, 1, $false, "30-10=20", 2)
[void]$d.Content.Find.Execute("3+57=60", $false, $false, $false, $false, $false, $true, 1, $false, "49+31=80", 2)
[void]$d.Content.Find.Execute("62-53=9", $false, $false, $false, $false, $false, $true, 1, $false, "99-23=76", 2)
[void]$d.Content.Find.Execute("67-2=65", $false, $false, $false, $false, $false, $true, 1, $false, "90-82=8", 2)
[void]$d.Content.Find.Execute("33-0=33", $false, $false, $false, $false, $false, $true, 1, $false, "64+23=87", 2)
[void]$d.Content.Find.Execute("85-23=62", $false, $false, $false, $false, $false, $true, 1, $false, "22+21=43", 2)
[void]$d.Content.Find.Execute("3+33=36", $false, $false, $false, $false, $false, $true, 1, $false, "72+15=87", 2)
[void]$d.Content.Find.Execute("58+37=95", $false, $false, $false, $false, $false, $true, 1, $false, "57-45=12", 2)
[void]$d.Content.Find.Execute("5+73=78", $false, $false, $false, $false, $false, $true, 1, $false, "1+23=24", 2)
[void]$d.Content.Find.Execute("18+64=82", $false, $false, $false, $false, $false, $true, 1, $false, "53+22=75", 2)
[void]$d.Content.Find.Execute("4+21=25", $false, $false, $false, $false, $false, $true, 1, $false, "6+67=73", 2)
[void]$d.Content.Find.Execute("59-5=54", $false, $false, $false, $false, $false, $true, 1, $false, "55+21=76", 2)
[void]$d.Content.Find.Execute("58-31=27", $false, $false, $false, $false, $false, $true, 1, $false, "4+42=46", 2)
[void]$d.Content.Find.Execute("55+25=80", $false, $false, $false, $false, $false, $true, 1, $false, "6+18=24", 2)
[void]$d.Content.Find.Execute("89+2=91", $false, $false, $false, $false, $false, $true, 1, $false, "11+28=39", 2)
[void]$d.Content.Find.Execute("97-60=37", $false, $false, $false, $false, $false, $true, 1, $false, "32+58=90", 2)
[void]$d.Content.Find.Execute("54-7=47", $false, $false, $false, $false, $false, $true, 1, $false, "68-4=64", 2)
[void]$d.Content.Find.Execute("17+62=79", $false, $false, $false, $false, $false, $true, 1, $false, "72-40=32", 2)
[void]$d.Content.Find.Execute("42+43=85", $false, $false, $false, $false, $false, $true, 1, $false, "45+14=59", 2)
[void]$d.Content.Find.Execute("91-14=77", $false, $false, $false, $false, $false, $true, 1, $false, "34+12=46", 2)
[void]$d.Content.Find.Execute("58+1=59", $false, $false, $false, $false, $false, $true, 1, $false, "43-2=41", 2)
[void]$d.Content.Find.Execute("46-17=29", $false, $false, $false, $false, $false, $true, 1, $false, "18+79=97", 2)
[void]$d.Content.Find.Execute("15-3=12", $false, $false, $false, $false, $false, $true, 1, $false, "9+57=66", 2)
[void]$d.Content.Find.Execute("86+0=86", $false, $false, $false, $false, $false, $true, 1, $false, "33-5=28", 2)
[void]$d.Content.Find.Execute("15-2=13", $false, $false, $false, $false, $false, $true, 1, $false, "40+40=80", 2)
[void]$d.Content.Find.Execute("77+8=85", $false, $false, $false, $false, $false, $true, 1, $false, "94-31=63", 2)
[void]$d.Content.Find.Execute("86-81=5", $false, $false, $false, $false, $false, $true, 1, $false, "11+61=72", 2)
[void]$d.Content.Find.Execute("59-13=46", $false, $false, $false, $false, $false, $true, 1, $false, "61+22=83", 2)
[void]$d.Content.Find.Execute("84-81=3", $false, $false, $false, $false, $false, $true, 1, $false, "72+20=92", 2)
[void]$d.Content.Find.Execute("48-47=1", $false, $false, $false, $false, $false, $true, 1, $false, "76-72=4", 2)
[void]$d.Content.Find.Execute("55+8=63", $false, $false, $false, $false, $false, $true, 1, $false, "54+23=77", 2)
[void]$d.Content.Find.Execute("83-52=31", $false, $false, $false, $false, $false, $true, 1, $false, "26+73=99", 2)
[void]$d.Content.Find.Execute("16+64=80", $false, $false, $false, $false, $false, $true, 1, $false, "25+59=84", 2)
[void]$d.Content.Find.Execute("44-8=36", $false, $false, $false, $false, $false, $true, 1, $false, "1+43=44", 2)
[void]$d.Content.Find.Execute("96-13=83", $false, $false, $false, $false, $false, $true, 1, $false, "94-14=80", 2)
[void]$d.Content.Find.Execute("1+54=55", $false, $false, $false, $false, $false, $true, 1, $false, "16+55=71", 2)
[void]$d.Content.Find.Execute("95-19=76", $false, $false, $false, $false, $false, $true, 1, $false, "65+24=89", 2)
[void]$d.Content.Find.Execute("5+44=49", $false, $false, $false, $false, $false, $true, 1, $false, "81-74=7", 2)
[void]$d.Content.Find.Execute("3+75=78", $false, $false, $false, $false, $false, $true, 1, $false, "79-36=43", 2)
[void]$d.Content.Find.Execute("40-33=7", $false, $false, $false, $false, $false, $true, 1, $false, "38-14=24", 2)
[void]$d.Content.Find.Execute("74-64=10", $false, $false, $false, $false, $false, $true, 1, $false, "47+50=97", 2)
[void]$d.Content.Find.Execute("15-5=10", $false, $false, $false, $false, $false, $true, 1, $false, "96-36=60", 2)
[void]$d.Content.Find.Execute("70-27=43", $false, $false, $false, $false, $false, $true, 1, $false, "35+29=64", 2)
[void]$d.Content.Find.Execute("68-40=28", $false, $false, $false, $false, $false, $true, 1, $false, "97-1=96", 2)
[void]$d.Content.Find.Execute("81-75=6", $false, $false, $false, $false, $false, $true, 1, $false, "0+16=16", 2)
[void]$d.Content.Find.Execute("10+39=49", $false, $false, $false, $false, $false, $true, 1, $false, "68-64=4", 2)
[void]$d.Content.Find.Execute("60+21=81", $false, $false, $false, $false, $false, $true, 1, $false, "89-53=36", 2)
[void]$d.Content.Find.Execute("14+42=56", $false, $false, $false, $false, $false, $true, 1, $false, "57-36=21", 2)
[void]$d.Content.Find.Execute("27+55=82", $false, $false, $false, $false, $false, $true, 1, $false, "11-1=10", 2)
[void]$d.Content.Find.Execute("84-50=34", $false, $false, $false, $false, $false, $true, 1, $false, "98-76=22", 2)
[void]$d.Content.Find.Execute("75+4=79", $false, $false, $false, $false, $false, $true, 1, $false, "57-16=41", 2)
[void]$d.Content.Find.Execute("58-6=52", $false, $false, $false, $false, $false, $true, 1, $false, "99-63=36", 2)
[void]$d.Content.Find.Execute("25+60=85", $false, $false, $false, $false, $false, $true, 1, $false, "81-73=8", 2)
[void]$d.Content.Find.Execute("90-40=50", $false, $false, $false, $false, $false, $true, 1, $false, "79-52=27", 2)
[void]$d.Content.Find.Execute("53-49=4", $false, $false, $false, $false, $false, $true, 1, $false, "3+80=83", 2)
[void]$d.Content.Find.Execute("12+55=67", $false, $false, $false, $false, $false, $true, 1, $false, "2+78=80", 2)
[void]$d.Content.Find.Execute("11+26=37", $false, $false, $false, $false, $false, $true, 1, $false, "67+16=83", 2)
[void]$d.Content.Find.Execute("16+31=47", $false, $false, $false, $false, $false, $true, 1, $false, "52-44=8", 2)
[void]$d.Content.Find.Execute("83-57=26", $false, $false, $false, $false, $false, $true, 1, $false, "11+43=54", 2)
[void]$d.Content.Find.Execute("76-19=57", $false, $false, $false, $false, $false, $true, 1, $false, "32+35=67", 2)
[void]$d.Content.Find.Execute("82+17=99", $false, $false, $false, $false, $false, $true, 1, $false, "33+14=47", 2)
[void]$d.Content.Find.Execute("50+0=50", $false, $false, $false, $false, $false, $true, 1, $false, "65-39=26", 2)
[void]$d.Content.Find.Execute("16+51=67", $false, $false, $false, $false, $false, $true, 1, $false, "22+74=96", 2)
[void]$d.Content.Find.Execute("8+34=42", $false, $false, $false, $false, $false, $true, 1, $false, "52+19=71", 2)
[void]$d.Content.Find.Execute("39-2=37", $false, $false, $false, $false, $false, $true, 1, $false, "25+57=82", 2)
[void]$d.Content.Find.Execute("66-52=14", $false, $false, $false, $false, $false, $true, 1, $false, "34-3=31", 2)
[void]$d.Content.Find.Execute("90-11=79", $false, $false, $false, $false, $false, $true, 1, $false, "80-41=39", 2)
[void]$d.Content.Find.Execute("18+3=21", $false, $false, $false, $false, $false, $true, 1, $false, "47-29=18", 2)
[void]$d.Content.Find.Execute("66+3=69", $false, $false, $false, $false, $false, $true, 1, $false, "73-67=6", 2)
[void]$d.Content.Find.Execute("41-9=32", $false, $false, $false, $false, $false, $true, 1, $false, "98-4=94", 2)
[void]$d.Content.Find.Execute("49-19=30", $false, $false, $false, $false, $false, $true, 1, $false, "12+62=74", 2)
[void]$d.Content.Find.Execute("80-44=36", $false, $false, $false, $false, $false, $true, 1, $false, "45-34=11", 2)
[void]$d.Content.Find.Execute("90-86=4", $false, $false, $false, $false, $false, $true, 1, $false, "25+72=97", 2)
[void]$d.Content.Find.Execute("63+8=71", $false, $false, $false, $false, $false, $true, 1, $false, "70-57=13", 2)
[void]$d.Content.Find.Execute("73-71=2", $false, $false, $false, $false, $false, $true, 1, $false, "81-74=7", 2)
[void]$d.Content.Find.Execute("65+30=95", $false, $false, $false, $false, $false, $true, 1, $false, "37+55=92", 2)
[void]$d.Content.Find.Execute("95-60=35", $false, $false, $false, $false, $false, $true, 1, $false, "1+26=27", 2)
[void]$d.Content.Find.Execute("18+71=89", $false, $false, $false, $false, $false, $true, 1, $false, "6+39=45", 2)
[void]$d.Content.Find.Execute("77-54=23", $false, $false, $false, $false, $false, $true, 1, $false, "44-27=17", 2)
[void]$d.Content.Find.Execute("95-43=52", $false, $false, $false, $false, $false, $true, 1, $false, "52-29=23", 2)
[void]$d.Content.Find.Execute("6+7=13", $false, $false, $false, $false, $false, $true, 1, $false, "39+37=76", 2)
[void]$d.Content.Find.Execute("79-74=5", $false, $false, $false, $false, $false, $true, 1, $false, "65-63=2", 2)
[void]$d.Content.Find.Execute("86-69=17", $false, $false, $false, $false, $false, $true, 1, $false, "93+6=99", 2)
[void]$d.Content.Find.Execute("69-5=64", $false, $false, $false, $false, $false, $true, 1, $false, "90-81=9", 2)
[void]$d.Content.Find.Execute("88+7=95", $false, $false, $false, $false, $false, $true, 1, $false, "30-0=30", 2)
[void]$d.Content.Find.Execute("96-89=7", $false, $false, $false, $false, $false, $true, 1, $false, "62+7=69", 2)
[void]$d.Content.Find.Execute("89-7=82", $false, $false, $false, $false, $false, $true, 1, $false, "16+42=58", 2)
[void]$d.Content.Find.Execute("32+65=97", $false, $false, $false, $false, $false, $true, 1, $false, "89-46=43", 2)
[void]$d.Content.Find.Execute("60+12=72", $false, $false, $false, $false, $false, $true, 1, $false, "32-24=8", 2)
[void]$d.Content.Find.Execute("52-19=33", $false, $false, $false, $false, $false, $true, 1, $false, "31+10=41", 2)
[void]$d.Content.Find.Execute("40+58=98", $false, $false, $false, $false, $false, $true, 1, $false, "9+35=44", 2)
[void]$d.Content.Find.Execute("93-81=12", $false, $false, $false, $false, $false, $true, 1, $false, "22+52=74", 2)
[void]$d.Content.Find.Execute("79-55=24", $false, $false, $false, $false, $false, $true, 1, $false, "46+30=76", 2)
[void]$d.Content.Find.Execute("66+13=79", $false, $false, $false, $false, $false, $true, 1, $false, "91-15=76", 2)
[void]$d.Content.Find.Execute("88-54=34", $false, $false, $false, $false, $false, $true, 1, $false, "21+5=26", 2)
[void]$d.Content.Find.Execute("19+19=38", $false, $false, $false, $false, $false, $true, 1, $false, "79+1=80", 2)
[void]$d.Content.Find.Execute("29+9=38", $false, $false, $false, $false, $false, $true, 1, $false, "33+32=65", 2)
[void]$d.Content.Find.Execute("62-50=12", $false, $false, $false, $false, $false, $true, 1, $false, "81-4=77", 2)
